$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1358.0834
$ws.Range("I19").Value = 1365.6666
$ws.Range("J19").Value = 1355.5555
$ws.Range("K19").Value = 1365.6666
$ws.Range("L19").Value = 1355.5555
$ws.Range("M19").Value = -1190.6666
$ws.Range("N19").Value = -1705.5555

# Row 32
$ws.Range("H32").Value = 2000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 2000
$ws.Range("N32").Value = -2652

# Row 39
$ws.Range("H39").Value = 603.35297
$ws.Range("I39").Value = 280.07693
$ws.Range("J39").Value = 1654
$ws.Range("K39").Value = 840.2307900000001
$ws.Range("L39").Value = 4962
$ws.Range("M39").Value = -544.2307900000001
$ws.Range("N39").Value = -5554

# Row 43
$ws.Range("H43").Value = 3488
$ws.Range("I43").Value = 1899.0769
$ws.Range("J43").Value = 7619.2
$ws.Range("K43").Value = 1899.0769
$ws.Range("L43").Value = 7619.2
$ws.Range("M43").Value = -1830.0769
$ws.Range("N43").Value = -7757.2

# Row 55
$ws.Range("H55").Value = 163
$ws.Range("I55").Value = 173.08333
$ws.Range("J55").Value = 138.8
$ws.Range("K55").Value = 173.08333
$ws.Range("L55").Value = 138.8
$ws.Range("M55").Value = 40.91667000000001
$ws.Range("N55").Value = -566.8

# Row 107
$ws.Range("H107").Value = 201.70589
$ws.Range("I107").Value = 201.70589
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 201.70589
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 1718.29411

# Row 137
$ws.Range("H137").Value = 4367.8696
$ws.Range("I137").Value = 1665.4667
$ws.Range("J137").Value = 9434.875
$ws.Range("K137").Value = 4996.4001
$ws.Range("L137").Value = 28304.625
$ws.Range("M137").Value = -2446.4001
$ws.Range("N137").Value = -33404.625

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2821.889
$ws.Range("I2").Value = 2560.625
$ws.Range("J2").Value = 4912
$ws.Range("K2").Value = 2560.625
$ws.Range("L2").Value = 4912
$ws.Range("M2").Value = -2447.625
$ws.Range("N2").Value = -5138

# Row 6
$ws.Range("H6").Value = 111543.445
$ws.Range("I6").Value = 143270.72
$ws.Range("J6").Value = 498
$ws.Range("K6").Value = 143270.72
$ws.Range("L6").Value = 498
$ws.Range("M6").Value = -143097.72
$ws.Range("N6").Value = -844

# Row 45
$ws.Range("H45").Value = 3700.7
$ws.Range("I45").Value = 3688.375
$ws.Range("J45").Value = 3750
$ws.Range("K45").Value = 3688.375
$ws.Range("L45").Value = 3750
$ws.Range("M45").Value = -3311.375
$ws.Range("N45").Value = -4504

# Row 61
$ws.Range("H61").Value = 20881260
$ws.Range("I61").Value = 26320304
$ws.Range("J61").Value = 212891.6
$ws.Range("K61").Value = 26320304
$ws.Range("L61").Value = 212891.6
$ws.Range("M61").Value = -26320092
$ws.Range("N61").Value = -213315.6

# Row 74
$ws.Range("H74").Value = 8937869
$ws.Range("I74").Value = 13891048
$ws.Range("J74").Value = 22146.8
$ws.Range("K74").Value = 13891048
$ws.Range("L74").Value = 22146.8
$ws.Range("M74").Value = -13890174
$ws.Range("N74").Value = -23894.8

# Row 77
$ws.Range("H77").Value = 8937869
$ws.Range("I77").Value = 13891048
$ws.Range("J77").Value = 22146.8
$ws.Range("K77").Value = 69455240
$ws.Range("L77").Value = 110734
$ws.Range("M77").Value = -69450872
$ws.Range("N77").Value = -119470

# Row 102
$ws.Range("H102").Value = 15411.143
$ws.Range("I102").Value = 17775.8
$ws.Range("J102").Value = 9499.5
$ws.Range("K102").Value = 17775.8
$ws.Range("L102").Value = 9499.5
$ws.Range("M102").Value = -16153.8
$ws.Range("N102").Value = -12743.5

# Row 110
$ws.Range("H110").Value = 3375
$ws.Range("I110").Value = 3239.5833
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 3239.5833
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = -1194.5833
$ws.Range("N110").Value = -9090

# Row 116
$ws.Range("H116").Value = 2821.889
$ws.Range("I116").Value = 2560.625
$ws.Range("J116").Value = 4912
$ws.Range("K116").Value = 2560.625
$ws.Range("L116").Value = 4912
$ws.Range("M116").Value = -266.625
$ws.Range("N116").Value = -9500

# Row 122
$ws.Range("H122").Value = 1311.2858
$ws.Range("I122").Value = 1408.5714
$ws.Range("J122").Value = 1214
$ws.Range("K122").Value = 4225.7142
$ws.Range("L122").Value = 3642
$ws.Range("M122").Value = -1775.7142
$ws.Range("N122").Value = -8542

# Row 132
$ws.Range("H132").Value = 6882.6
$ws.Range("I132").Value = 4798.7393
$ws.Range("J132").Value = 13729.571
$ws.Range("K132").Value = 14396.2179
$ws.Range("L132").Value = 41188.713
$ws.Range("M132").Value = -11866.2179
$ws.Range("N132").Value = -46248.713

# Row 136
$ws.Range("H136").Value = 20881260
$ws.Range("I136").Value = 26320304
$ws.Range("J136").Value = 212891.6
$ws.Range("K136").Value = 78960912
$ws.Range("L136").Value = 638674.8
$ws.Range("M136").Value = -78958362
$ws.Range("N136").Value = -643774.8

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2821.889
$ws.Range("I3").Value = 2560.625
$ws.Range("J3").Value = 4912
$ws.Range("K3").Value = 2560.625
$ws.Range("L3").Value = 4912
$ws.Range("M3").Value = -2446.625
$ws.Range("N3").Value = -5140

# Row 20
$ws.Range("H20").Value = 3019.4375
$ws.Range("I20").Value = 2732.6924
$ws.Range("J20").Value = 4262
$ws.Range("K20").Value = 2732.6924
$ws.Range("L20").Value = 4262
$ws.Range("M20").Value = -2485.6924
$ws.Range("N20").Value = -4756

# Row 59
$ws.Range("H59").Value = 99999
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 99999
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 99999
$ws.Range("N59").Value = -101693

# Row 107
$ws.Range("H107").Value = 838.1905
$ws.Range("I107").Value = 805.2
$ws.Range("J107").Value = 1498
$ws.Range("K107").Value = 805.2
$ws.Range("L107").Value = 1498
$ws.Range("M107").Value = 1114.8
$ws.Range("N107").Value = -5338

# Row 134
$ws.Range("H134").Value = 33871.39
$ws.Range("I134").Value = 2249.85
$ws.Range("J134").Value = 73398.31
$ws.Range("K134").Value = 6749.549999999999
$ws.Range("L134").Value = 220194.93
$ws.Range("M134").Value = -4214.549999999999
$ws.Range("N134").Value = -225264.93

$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Range("H12").Value = 50
$ws.Range("I12").Value = 50
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 50
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 120

# Row 31
$ws.Range("H31").Value = 443695.16
$ws.Range("I31").Value = 3295.76
$ws.Range("J31").Value = 836908.9
$ws.Range("K31").Value = 3295.76
$ws.Range("L31").Value = 836908.9
$ws.Range("M31").Value = -3000.76
$ws.Range("N31").Value = -837498.9

# Row 34
$ws.Range("H34").Value = 443695.16
$ws.Range("I34").Value = 3295.76
$ws.Range("J34").Value = 836908.9
$ws.Range("K34").Value = 3295.76
$ws.Range("L34").Value = 836908.9
$ws.Range("M34").Value = -3093.76
$ws.Range("N34").Value = -837312.9

# Row 92
$ws.Range("H92").Value = 48330
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 48330
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 48330
$ws.Range("N92").Value = -53322

# Row 134
$ws.Range("H134").Value = 297562.2
$ws.Range("I134").Value = 386347.38
$ws.Range("J134").Value = 9010.25
$ws.Range("K134").Value = 1159042.14
$ws.Range("L134").Value = 27030.75
$ws.Range("M134").Value = -1156507.14
$ws.Range("N134").Value = -32100.75

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 2013.88
$ws.Range("I11").Value = 2047.7916
$ws.Range("J11").Value = 1200
$ws.Range("K11").Value = 6143.3748
$ws.Range("L11").Value = 3600
$ws.Range("M11").Value = -6003.3748
$ws.Range("N11").Value = -3880

# Row 103
$ws.Range("H103").Value = 1363.125
$ws.Range("I103").Value = 474
$ws.Range("J103").Value = 1896.6
$ws.Range("K103").Value = 1422
$ws.Range("L103").Value = 5689.799999999999
$ws.Range("M103").Value = -543
$ws.Range("N103").Value = -7447.799999999999

# Row 113
$ws.Range("H113").Value = 2240.1667
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 2860.25
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 8580.75
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -12920.75

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 6718.5
$ws.Range("I102").Value = 4106.857
$ws.Range("J102").Value = 25000
$ws.Range("K102").Value = 4106.857
$ws.Range("L102").Value = 25000
$ws.Range("M102").Value = -2484.857
$ws.Range("N102").Value = -28244

# Row 122
$ws.Range("H122").Value = 2226.6667
$ws.Range("I122").Value = 2446.75
$ws.Range("J122").Value = 1786.5
$ws.Range("K122").Value = 7340.25
$ws.Range("L122").Value = 5359.5
$ws.Range("M122").Value = -4890.25
$ws.Range("N122").Value = -10259.5

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2262.5454
$ws.Range("I22").Value = 2411.875
$ws.Range("J22").Value = 1864.3334
$ws.Range("K22").Value = 2411.875
$ws.Range("L22").Value = 1864.3334
$ws.Range("M22").Value = -2116.875
$ws.Range("N22").Value = -2454.3334

# Row 27
$ws.Range("H27").Value = 2262.5454
$ws.Range("I27").Value = 2411.875
$ws.Range("J27").Value = 1864.3334
$ws.Range("K27").Value = 2411.875
$ws.Range("L27").Value = 1864.3334
$ws.Range("M27").Value = -2304.875
$ws.Range("N27").Value = -2078.3334

# Row 36
$ws.Range("H36").Value = 77280.8
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 77280.8
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 77280.8
$ws.Range("N36").Value = -78404.8

# Row 61
$ws.Range("H61").Value = 631
$ws.Range("I61").Value = 631
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 631
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -429

# Row 93
$ws.Range("H93").Value = 47620296
$ws.Range("I93").Value = 52632732
$ws.Range("J93").Value = 2150.5
$ws.Range("K93").Value = 52632732
$ws.Range("L93").Value = 2150.5
$ws.Range("M93").Value = -52631484
$ws.Range("N93").Value = -4646.5

# Row 100
$ws.Range("H100").Value = 4001.6667
$ws.Range("I100").Value = 4001.5
$ws.Range("J100").Value = 4002
$ws.Range("K100").Value = 4001.5
$ws.Range("L100").Value = 4002
$ws.Range("M100").Value = -3460.5
$ws.Range("N100").Value = -5084

# Row 113
$ws.Range("H113").Value = 631
$ws.Range("I113").Value = 631
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 631
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1539

# Row 122
$ws.Range("H122").Value = 5649.3887
$ws.Range("I122").Value = 4399.154
$ws.Range("J122").Value = 8900
$ws.Range("K122").Value = 13197.462
$ws.Range("L122").Value = 26700
$ws.Range("M122").Value = -10747.462
$ws.Range("N122").Value = -31600

# Row 136
$ws.Range("H136").Value = 51200.24
$ws.Range("I136").Value = 5759.7646
$ws.Range("J136").Value = 147761.25
$ws.Range("K136").Value = 17279.2938
$ws.Range("L136").Value = 443283.75
$ws.Range("M136").Value = -14729.2938
$ws.Range("N136").Value = -448383.75

# Row 140
$ws.Range("H140").Value = 111792.336
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 111792.336
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 111792.336
$ws.Range("N140").Value = -122152.336

$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 26666.334
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 26666.334
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 26666.334
$ws.Range("N39").Value = -27492.334

# Row 43
$ws.Range("H43").Value = 80000
$ws.Range("I43").Value = 80000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 80000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -79851

# Row 96
$ws.Range("H96").Value = 3662.125
$ws.Range("I96").Value = 2649.75
$ws.Range("J96").Value = 4674.5
$ws.Range("K96").Value = 2649.75
$ws.Range("L96").Value = 4674.5
$ws.Range("M96").Value = -1276.75
$ws.Range("N96").Value = -7420.5

# Row 113
$ws.Range("H113").Value = 699.13043
$ws.Range("I113").Value = 744.5714
$ws.Range("J113").Value = 222
$ws.Range("K113").Value = 2233.7142
$ws.Range("L113").Value = 666
$ws.Range("M113").Value = -63.71420000000035
$ws.Range("N113").Value = -5006

# Row 122
$ws.Range("H122").Value = 6464.4116
$ws.Range("I122").Value = 1984.6
$ws.Range("J122").Value = 12864.143
$ws.Range("K122").Value = 5953.799999999999
$ws.Range("L122").Value = 38592.429
$ws.Range("M122").Value = -3503.799999999999
$ws.Range("N122").Value = -43492.429

# Row 132
$ws.Range("H132").Value = 7669.1724
$ws.Range("I132").Value = 961.6087
$ws.Range("J132").Value = 33381.5
$ws.Range("K132").Value = 2884.8261
$ws.Range("L132").Value = 100144.5
$ws.Range("M132").Value = -354.8261000000002
$ws.Range("N132").Value = -105204.5
